$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (mmr) row 2 values before overwriting
$row2Values = @()
for ($col = 2; $col -le 13; $col++) {
    $row2Values += $ws.Cells.Item(2, $col).Value
}

# Row 2 becomes "mf": label + D:M values from old row 3, B and C blank
$ws.Cells.Item(2, 1).Value = "mf"
$ws.Cells.Item(2, 2).Value = ""
$ws.Cells.Item(2, 3).Value = ""
$ws.Cells.Item(2, 4).Value = 0.03011123897199847
$ws.Cells.Item(2, 5).Value = 0.180425360411634
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 0.03011123897199847

# Row 3 becomes "mmr": label + values previously in row 2
$ws.Cells.Item(3, 1).Value = "mmr"
$ws.Cells.Item(3, 2).Value = 0.7134361190331808
$ws.Cells.Item(3, 3).Value = 0.5110410933649582
$ws.Cells.Item(3, 4).Value = 0.03778289221327196
$ws.Cells.Item(3, 5).Value = 0.1931026812795729
$ws.Cells.Item(3, 6).Value = 0.00909090909090909
$ws.Cells.Item(3, 7).Value = 0.005726110502104429
$ws.Cells.Item(3, 8).Value = 0.005227462901366726
$ws.Cells.Item(3, 9).Value = 0.008377837609220288
$ws.Cells.Item(3, 10).Value = 0.00231934626038035
$ws.Cells.Item(3, 11).Value = 0.01796536796536797
$ws.Cells.Item(3, 12).Value = 0.003506940450404207
$ws.Cells.Item(3, 13).Value = 0.03778289221327196
